$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test Case blocks are laid out as:
#   TC2 block: header row 15, step rows 19-20
#   TC3 block: header row 23, step rows 27-28
#   TC4 block: header row 31, step rows 35-36
#
# The edit rotates the 2nd step/expected pair among the TC2, TC3 and TC4
# blocks so that:
#   TC2's 2nd step becomes "ordenar pelo nome do servidor" (was TC4's step)
#   TC3's 2nd step becomes "filtrar a listagem"             (was TC2's step)
#   TC4's 2nd step becomes "cancelamento de uma diaria"     (was TC3's step)

$ws.Range("B20").Value = "Chefe Clica para ordenar pelo nome do servidor."
$ws.Range("D20").Value = "SYSTEM Visualiza os registros de solicitações de diária ordenado pelo nome do servidor."

$ws.Range("B28").Value = "Chefe Indica alguns parâmetros específicos para a busca; Informa o nome do beneficiário; Filtra a listagem de solicitações."
$ws.Range("D28").Value = "SYSTEM Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário."

$ws.Range("B36").Value = "Chefe Clica para realizar o cancelamento de uma diária."
$ws.Range("D36").Value = "SYSTEM Verifica que a solicitação está em situação SOLICITADA; Exibe mensagem de confirmação (MSG987 - Cancelar solicitação de diária) para o usuário (que deve confirmar); Cancela a diária, mudando sua situação para CANCELADA (ver diagrama de estados da diária)."
